$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Name: "Malc" + _GoBack bookmark + "olm Wanstall"  ->  "Malcolm Wanstall"
#    (single merged run, bookmark disappears from here - it gets re-added
#    later near the Medibank employer line)
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Malcolm Wanstall", $true, $false, $false, $false, $false, $true, 1, $false, "Malcolm Wanstall", 2)

# ---------------------------------------------------------------------------
# 2) Summary paragraph: merge the three runs ("...already ", "possess",
#    " but have yet to fully leverage.") into a single run.
# ---------------------------------------------------------------------------
$summarySentence = "I thrive on turning vision into reality, to help companies realise the hidden assets they have in the data, people and capabilities they often already possess but have yet to fully leverage."
$null = $d.Content.Find.Execute($summarySentence, $true, $false, $false, $false, $false, $true, 1, $false, $summarySentence, 2)

# ---------------------------------------------------------------------------
# 3) Job title: "Data Warehouse Manager / EDW/BI Product Owner"
#    -> "Data Warehouse " + "/ Business Intelligence Program Manager"
#    (two runs, identical rPr) - use TrackRevisions so the engine keeps the
#    insertion as a distinct run instead of silently re-merging it with its
#    neighbour.
# ---------------------------------------------------------------------------
$was_tracking = $d.TrackRevisions
$d.TrackRevisions = $true

$rng = $d.Content
$null = $rng.Find.Execute("Data Warehouse Manager / EDW/BI Product Owner", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$matchStart = $rng.Start

$insertPoint = $d.Range($matchStart + 15, $matchStart + 15)
$insertPoint.InsertAfter("/ Business Intelligence Program Manager")

$tailRng = $d.Content
$null = $tailRng.Find.Execute("Manager / EDW/BI Product Owner", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailRng.Text = ""

# ---------------------------------------------------------------------------
# 4) Employer line: "MEDIBANK HEALTH SOLUTIONS - NOV 2011 -> PRESENT"
#    -> "MEDIBANK " + "(MHS)" + [_GoBack bookmark] + " - NOV 2011 -> PRESENT"
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$null = $rng2.Find.Execute("MEDIBANK HEALTH SOLUTIONS - NOV 2011", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$empStart = $rng2.Start

$insertPoint2 = $d.Range($empStart + 9, $empStart + 9)
$insertPoint2.InsertAfter("(MHS)")

$tailRng2 = $d.Content
$null = $tailRng2.Find.Execute("HEALTH SOLUTIONS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailRng2.Text = ""

$d.TrackRevisions = $was_tracking
$d.Revisions.AcceptAll()

# Re-add the _GoBack bookmark between "(MHS)" and " - NOV 2011" on the
# employer line.
$bmRng = $d.Content
$null = $bmRng.Find.Execute("(MHS)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPoint = $d.Range($bmRng.End, $bmRng.End)
$null = $d.Bookmarks.Add("_GoBack", $bmPoint)

Write-Output "done"
